# Update NATMI LR-pair sheet with newly recomputed TPM values.
# Adds "ECs" as a new sending cluster (rows 2-4) and shifts the former
# Inflammatory-Mac / Neutrophils / Resolving-Mac blocks down, appending a
# new Resolving-Mac block (rows 11-13) that did not exist before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 1.005755333333333
$ws.Range("H2").Value = 3.017266
$ws.Range("I2").Value = 0.01048729000197281
$ws.Range("J2").Value = 0.01048729000197281
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05145833333333334
$ws.Range("N2").Value = 0.154375
$ws.Range("O2").Value = 0.2409462730781657
$ws.Range("P2").Value = 0.2409462730781657
$ws.Range("Q2").Value = 0.05175449319444445
$ws.Range("R2").Value = 0.46579043875
$ws.Range("S2").Value = 0.002526873440665257
$ws.Range("T2").Value = 0.002526873440665257

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 1.005755333333333
$ws.Range("H3").Value = 3.017266
$ws.Range("I3").Value = 0.01048729000197281
$ws.Range("J3").Value = 0.01048729000197281
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1394176666666667
$ws.Range("N3").Value = 0.418253
$ws.Range("O3").Value = 0.6528032489312521
$ws.Range("P3").Value = 0.6528032489312521
$ws.Range("Q3").Value = 0.1402200618108889
$ws.Range("R3").Value = 1.261980556298
$ws.Range("S3").Value = 0.006846136985772085
$ws.Range("T3").Value = 0.006846136985772085

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ccl12"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 1.005755333333333
$ws.Range("H4").Value = 3.017266
$ws.Range("I4").Value = 0.01048729000197281
$ws.Range("J4").Value = 0.01048729000197281
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02269166666666667
$ws.Range("N4").Value = 0.068075
$ws.Range("O4").Value = 0.1062504779905822
$ws.Range("P4").Value = 0.1062504779905822
$ws.Range("Q4").Value = 0.02282226477222222
$ws.Range("R4").Value = 0.20540038295
$ws.Range("S4").Value = 0.001114279575535464
$ws.Range("T4").Value = 0.001114279575535464

# Row 5
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Ccl12"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 38.755371
$ws.Range("H5").Value = 116.266113
$ws.Range("I5").Value = 0.4041130097356814
$ws.Range("J5").Value = 0.4041130097356814
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05145833333333334
$ws.Range("N5").Value = 0.154375
$ws.Range("O5").Value = 0.2409462730781657
$ws.Range("P5").Value = 0.2409462730781657
$ws.Range("Q5").Value = 1.994286799375
$ws.Range("R5").Value = 17.948581194375
$ws.Range("S5").Value = 0.09736952359821292
$ws.Range("T5").Value = 0.09736952359821292

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Ccl12"
$ws.Range("C6").Value = "Ccr10"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 38.755371
$ws.Range("H6").Value = 116.266113
$ws.Range("I6").Value = 0.4041130097356814
$ws.Range("J6").Value = 0.4041130097356814
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1394176666666667
$ws.Range("N6").Value = 0.418253
$ws.Range("O6").Value = 0.6528032489312521
$ws.Range("P6").Value = 0.6528032489312521
$ws.Range("Q6").Value = 5.403183395620999
$ws.Range("R6").Value = 48.62865056058899
$ws.Range("S6").Value = 0.2638062856908395
$ws.Range("T6").Value = 0.2638062856908395

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Ccl12"
$ws.Range("C7").Value = "Ccr10"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 38.755371
$ws.Range("H7").Value = 116.266113
$ws.Range("I7").Value = 0.4041130097356814
$ws.Range("J7").Value = 0.4041130097356814
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02269166666666667
$ws.Range("N7").Value = 0.068075
$ws.Range("O7").Value = 0.1062504779905822
$ws.Range("P7").Value = 0.1062504779905822
$ws.Range("Q7").Value = 0.8794239602749999
$ws.Range("R7").Value = 7.914815642474999
$ws.Range("S7").Value = 0.04293720044662895
$ws.Range("T7").Value = 0.04293720044662895

# Row 8
$ws.Range("A8").Value = "Neutrophils"
$ws.Range("B8").Value = "Ccl12"
$ws.Range("C8").Value = "Ccr10"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.988471999999999
$ws.Range("H8").Value = 23.965416
$ws.Range("I8").Value = 0.08329801469605898
$ws.Range("J8").Value = 0.08329801469605898
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.05145833333333334
$ws.Range("N8").Value = 0.154375
$ws.Range("O8").Value = 0.2409462730781657
$ws.Range("P8").Value = 0.2409462730781657
$ws.Range("Q8").Value = 0.411073455
$ws.Range("R8").Value = 3.699661095
$ws.Range("S8").Value = 0.02007034619582569
$ws.Range("T8").Value = 0.02007034619582569

# Row 9
$ws.Range("A9").Value = "Neutrophils"
$ws.Range("B9").Value = "Ccl12"
$ws.Range("C9").Value = "Ccr10"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.988471999999999
$ws.Range("H9").Value = 23.965416
$ws.Range("I9").Value = 0.08329801469605898
$ws.Range("J9").Value = 0.08329801469605898
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1394176666666667
$ws.Range("N9").Value = 0.418253
$ws.Range("O9").Value = 0.6528032489312521
$ws.Range("P9").Value = 0.6528032489312521
$ws.Range("Q9").Value = 1.113734126472
$ws.Range("R9").Value = 10.023607138248
$ws.Range("S9").Value = 0.05437721462311048
$ws.Range("T9").Value = 0.05437721462311048

# Row 10
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "Ccl12"
$ws.Range("C10").Value = "Ccr10"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.988471999999999
$ws.Range("H10").Value = 23.965416
$ws.Range("I10").Value = 0.08329801469605898
$ws.Range("J10").Value = 0.08329801469605898
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02269166666666667
$ws.Range("N10").Value = 0.068075
$ws.Range("O10").Value = 0.1062504779905822
$ws.Range("P10").Value = 0.1062504779905822
$ws.Range("Q10").Value = 0.1812717438
$ws.Range("R10").Value = 1.6314456942
$ws.Range("S10").Value = 0.00885045387712281
$ws.Range("T10").Value = 0.008850453877122808

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Ccl12"
$ws.Range("C11").Value = "Ccr10"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 48.15271133333334
$ws.Range("H11").Value = 144.458134
$ws.Range("I11").Value = 0.5021016855662869
$ws.Range("J11").Value = 0.5021016855662868
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.05145833333333334
$ws.Range("N11").Value = 0.154375
$ws.Range("O11").Value = 0.2409462730781657
$ws.Range("P11").Value = 0.2409462730781657
$ws.Range("Q11").Value = 2.477858270694445
$ws.Range("R11").Value = 22.30072443625
$ws.Range("S11").Value = 0.1209795298434619
$ws.Range("T11").Value = 0.1209795298434618

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Ccl12"
$ws.Range("C12").Value = "Ccr10"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 48.15271133333334
$ws.Range("H12").Value = 144.458134
$ws.Range("I12").Value = 0.5021016855662869
$ws.Range("J12").Value = 0.5021016855662868
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.1394176666666667
$ws.Range("N12").Value = 0.418253
$ws.Range("O12").Value = 0.6528032489312521
$ws.Range("P12").Value = 0.6528032489312521
$ws.Range("Q12").Value = 6.713338657766889
$ws.Range("R12").Value = 60.420047919902
$ws.Range("S12").Value = 0.32777361163153
$ws.Range("T12").Value = 0.32777361163153

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Ccl12"
$ws.Range("C13").Value = "Ccr10"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 48.15271133333334
$ws.Range("H13").Value = 144.458134
$ws.Range("I13").Value = 0.5021016855662869
$ws.Range("J13").Value = 0.5021016855662868
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.02269166666666667
$ws.Range("N13").Value = 0.068075
$ws.Range("O13").Value = 0.1062504779905822
$ws.Range("P13").Value = 0.1062504779905822
$ws.Range("Q13").Value = 1.092665274672222
$ws.Range("R13").Value = 9.83398747205
$ws.Range("S13").Value = 0.053348544091295
$ws.Range("T13").Value = 0.05334854409129498
